$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 410
$ws.Range("I11").Value = 410
$ws.Range("K11").Value = 410
$ws.Range("M11").Value = -270
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H28").Value = 10287.7
$ws.Range("I28").Value = 14495
$ws.Range("J28").Value = 470.66666
$ws.Range("K28").Value = 14495
$ws.Range("L28").Value = 470.66666
$ws.Range("M28").Value = -14010
$ws.Range("N28").Value = -1440.66666
$ws.Range("H33").Value = 358.0476
$ws.Range("I33").Value = 360.55554
$ws.Range("K33").Value = 360.55554
$ws.Range("M33").Value = -131.55554
$ws.Range("H76").Value = 3518.5
$ws.Range("J76").Value = 3600
$ws.Range("L76").Value = 3600
$ws.Range("N76").Value = -4230
$ws.Range("H79").Value = 3518.5
$ws.Range("J79").Value = 3600
$ws.Range("L79").Value = 3600
$ws.Range("N79").Value = -5784
$ws.Range("H112").Value = 2032.1041
$ws.Range("J112").Value = 2072.6304
$ws.Range("L112").Value = 6217.8912
$ws.Range("N112").Value = -8433.8912
$ws.Range("H113").Value = 2093.2856
$ws.Range("I113").Value = 1775.5
$ws.Range("K113").Value = 1775.5
$ws.Range("M113").Value = 1478.5
$ws.Range("H132").Value = 8778933
$ws.Range("I132").Value = 14495949
$ws.Range("J132").Value = 12841.2
$ws.Range("K132").Value = 43487847
$ws.Range("L132").Value = 38523.60000000001
$ws.Range("M132").Value = -43485317
$ws.Range("N132").Value = -43583.60000000001
$ws.Range("H137").Value = 1067.3334
$ws.Range("I137").Value = 999
$ws.Range("K137").Value = 2997
$ws.Range("M137").Value = -447
$ws.Range("H138").Value = 751508.1
$ws.Range("J138").Value = 990211.3
$ws.Range("L138").Value = 2970633.9
$ws.Range("N138").Value = -2980913.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1529.7391
$ws.Range("I61").Value = 1425.4736
$ws.Range("J61").Value = 2025
$ws.Range("K61").Value = 1425.4736
$ws.Range("L61").Value = 2025
$ws.Range("M61").Value = -1213.4736
$ws.Range("N61").Value = -2449
$ws.Range("H88").Value = 2518.6667
$ws.Range("I88").Value = 2082.4
$ws.Range("J88").Value = 2736.8
$ws.Range("K88").Value = 2082.4
$ws.Range("L88").Value = 2736.8
$ws.Range("M88").Value = -1676.4
$ws.Range("N88").Value = -3548.8
$ws.Range("H91").Value = 2518.6667
$ws.Range("I91").Value = 2082.4
$ws.Range("J91").Value = 2736.8
$ws.Range("K91").Value = 2082.4
$ws.Range("L91").Value = 2736.8
$ws.Range("M91").Value = -678.4000000000001
$ws.Range("N91").Value = -5544.8
$ws.Range("H122").Value = 1250.2142
$ws.Range("I122").Value = 1224.3
$ws.Range("J122").Value = 1315
$ws.Range("K122").Value = 3672.9
$ws.Range("L122").Value = 3945
$ws.Range("M122").Value = -1222.9
$ws.Range("N122").Value = -8845
$ws.Range("H136").Value = 1529.7391
$ws.Range("I136").Value = 1425.4736
$ws.Range("J136").Value = 2025
$ws.Range("K136").Value = 4276.4208
$ws.Range("L136").Value = 6075
$ws.Range("M136").Value = -1726.4208
$ws.Range("N136").Value = -11175

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 50001600
$ws.Range("I99").Value = 50001600
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 50001600
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -50000102
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 440.22223
$ws.Range("I7").Value = 485.25
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 485.25
$ws.Range("L7").Value = 80
$ws.Range("M7").Value = -372.25
$ws.Range("N7").Value = -306
$ws.Range("H22").Value = 399
$ws.Range("I22").Value = 398.9
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 398.9
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -48.89999999999998
$ws.Range("N22").Value = -1100
$ws.Range("H31").Value = 1565.0625
$ws.Range("I31").Value = 1430.5
$ws.Range("J31").Value = 2507
$ws.Range("K31").Value = 1430.5
$ws.Range("L31").Value = 2507
$ws.Range("M31").Value = -1135.5
$ws.Range("N31").Value = -3097
$ws.Range("H34").Value = 1565.0625
$ws.Range("I34").Value = 1430.5
$ws.Range("J34").Value = 2507
$ws.Range("K34").Value = 1430.5
$ws.Range("L34").Value = 2507
$ws.Range("M34").Value = -1228.5
$ws.Range("N34").Value = -2911
$ws.Range("H86").Value = 6082517.5
$ws.Range("I86").Value = 9527610
$ws.Range("K86").Value = 9527610
$ws.Range("M86").Value = -9526487
$ws.Range("H89").Value = 6082517.5
$ws.Range("I89").Value = 9527610
$ws.Range("K89").Value = 47638050
$ws.Range("M89").Value = -47632434
$ws.Range("H99").Value = 1671.9166
$ws.Range("I99").Value = 1666.3
$ws.Range("J99").Value = 1700
$ws.Range("K99").Value = 1666.3
$ws.Range("L99").Value = 1700
$ws.Range("M99").Value = -168.3
$ws.Range("N99").Value = -4696
$ws.Range("H122").Value = 1002.7143
$ws.Range("I122").Value = 984.4
$ws.Range("J122").Value = 1048.5
$ws.Range("K122").Value = 2953.2
$ws.Range("L122").Value = 3145.5
$ws.Range("M122").Value = -503.1999999999998
$ws.Range("N122").Value = -8045.5
$ws.Range("H126").Value = 1671.9166
$ws.Range("I126").Value = 1666.3
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 4998.9
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -2528.9
$ws.Range("N126").Value = -10040
$ws.Range("H132").Value = 8140.3887
$ws.Range("I132").Value = 11486.6
$ws.Range("J132").Value = 3957.625
$ws.Range("K132").Value = 34459.8
$ws.Range("L132").Value = 11872.875
$ws.Range("M132").Value = -31929.8
$ws.Range("N132").Value = -16932.875
$ws.Range("H134").Value = 13334846
$ws.Range("I134").Value = 15152916
$ws.Range("K134").Value = 45458748
$ws.Range("M134").Value = -45456213

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 6516.25
$ws.Range("I120").Value = 999.5
$ws.Range("K120").Value = 2998.5
$ws.Range("M120").Value = 1839.5
$ws.Range("H122").Value = 732.0769
$ws.Range("I122").Value = 695.5714
$ws.Range("J122").Value = 774.6667
$ws.Range("K122").Value = 6260.1426
$ws.Range("L122").Value = 6972.0003
$ws.Range("M122").Value = -3810.1426
$ws.Range("N122").Value = -11872.0003
$ws.Range("H131").Value = 19232134
$ws.Range("J131").Value = 1572.6976
$ws.Range("L131").Value = 4718.0928
$ws.Range("N131").Value = -14798.0928
$ws.Range("H134").Value = 4775.6665
$ws.Range("I134").Value = 2145.8
$ws.Range("J134").Value = 5682.517
$ws.Range("K134").Value = 6437.400000000001
$ws.Range("L134").Value = 17047.551
$ws.Range("M134").Value = -1367.400000000001
$ws.Range("N134").Value = -27187.551

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 4273852
$ws.Range("I12").Value = 4607454.5
$ws.Range("J12").Value = 2806000
$ws.Range("K12").Value = 4607454.5
$ws.Range("L12").Value = 2806000
$ws.Range("M12").Value = -4607314.5
$ws.Range("N12").Value = -2806280
$ws.Range("H70").Value = 15004677
$ws.Range("I70").Value = 13162037
$ws.Range("J70").Value = 18187418
$ws.Range("K70").Value = 13162037
$ws.Range("L70").Value = 18187418
$ws.Range("M70").Value = -13161767
$ws.Range("N70").Value = -18187958
$ws.Range("H73").Value = 15004677
$ws.Range("I73").Value = 13162037
$ws.Range("J73").Value = 18187418
$ws.Range("K73").Value = 13162037
$ws.Range("L73").Value = 18187418
$ws.Range("M73").Value = -13161101
$ws.Range("N73").Value = -18189290
$ws.Range("H80").Value = 4568.8667
$ws.Range("I80").Value = 2100
$ws.Range("J80").Value = 6214.778
$ws.Range("K80").Value = 2100
$ws.Range("L80").Value = 6214.778
$ws.Range("M80").Value = -1102
$ws.Range("N80").Value = -8210.778
$ws.Range("H83").Value = 4568.8667
$ws.Range("I83").Value = 2100
$ws.Range("J83").Value = 6214.778
$ws.Range("K83").Value = 10500
$ws.Range("L83").Value = 31073.89
$ws.Range("M83").Value = -5508
$ws.Range("N83").Value = -41057.89

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 171815.17
$ws.Range("I132").Value = 6444.5
$ws.Range("J132").Value = 254500.5
$ws.Range("K132").Value = 19333.5
$ws.Range("L132").Value = 763501.5
$ws.Range("M132").Value = -16803.5
$ws.Range("N132").Value = -768561.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 10814.462
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 10814.462
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 10814.462
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -12062.462
$ws.Range("H66").Value = 10814.462
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 10814.462
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 32443.386
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -38683.386
$ws.Range("H96").Value = 2130
$ws.Range("J96").Value = 2080
$ws.Range("L96").Value = 2080
$ws.Range("N96").Value = -4826
$ws.Range("H132").Value = 4847.4707
$ws.Range("I132").Value = 5040.6
$ws.Range("J132").Value = 4571.5713
$ws.Range("K132").Value = 15121.8
$ws.Range("L132").Value = 13714.7139
$ws.Range("M132").Value = -12591.8
$ws.Range("N132").Value = -18774.7139
$ws.Range("H136").Value = 703.0769
$ws.Range("I136").Value = 304.44446
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 913.33338
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = 1636.66662
$ws.Range("N136").Value = -9900
